$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (column G) value
$kValues = @{
    2  = 4
    3  = 6
    4  = 3
    5  = 8
    6  = 2
    7  = 3
    8  = 3
    9  = 6
    10 = 5
    11 = 6
    12 = 1
    13 = 2
    14 = 3
    15 = 4
    16 = 6
    17 = 2
    18 = 3
    19 = 3
    20 = 4
    21 = 4
    22 = 5
    23 = 3
    24 = 0
    25 = 1
    26 = 0
    27 = 3
    28 = 0
    29 = 1
    30 = 1
    31 = 6
    32 = 5
    33 = 4
    34 = 3
    35 = 3
    36 = 4
    37 = 7
    38 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
